# Test Data generation Konzept
# Reworks the "Daten" worksheet: adds a Phi/Delta block above the existing
# Radius/angle table, and extends each angle row with direction-vector and
# start-point columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert two new rows -------------------------------------
# Row 2 becomes a new "Phi" row (old row 2 "# / Radius / 10" shifts to row 3).
# Row 4 becomes a new "Delta" row (old header row 3 shifts to row 5, and the
# angle data that used to start at row 4 now starts at row 6).
$ws.Rows("2:2").Insert()
$ws.Rows("4:4").Insert()

# --- 2. Header / parameter block (rows 1-5) ---------------------------------
$ws.Range("A2").Value = "#"
$ws.Range("C2").Value = 89.99

$ws.Range("A4").Value = "x"
$ws.Range("C4").Formula = "=3*C3"

# --- 3. Data rows (rows 6-14) ------------------------------------------------
$angles = 0, 45, 90, 135, 180, 225, 270, 315, 360

for ($i = 0; $i -lt $angles.Length; $i++) {
    $r = 6 + $i
    $ws.Cells.Item($r, 1).Value = $angles[$i]
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).NumberFormat = "0.0"
    $ws.Cells.Item($r, 3).Formula = "=`$C`$3*SIN(A$r*PI()/180)"
    $ws.Cells.Item($r, 4).NumberFormat = "0.0"
    $ws.Cells.Item($r, 4).Formula = "=`$C`$3*COS(A$r*PI()/180)"
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
}

# Extra column only present on the first data row: distance from the tip of
# the Delta segment to the point on the circle at angle Phi.
$ws.Range("H6").Formula = '=$C$4-$C$3*COS(PI()*$C$2/180)'

# --- 4. Extend the header row (row 5) with the new vector/point columns ----
# Written in this particular order so newly-introduced labels land in the
# same shared-string slots as in the authored workbook.
$ws.Range("E5").Value = "xDir"
$ws.Range("B2").Value = "Phi"
$ws.Range("F5").Value = "yDir"
$ws.Range("G5").Value = "zDir"
$ws.Range("H5").Value = "xPoint"
$ws.Range("I5").Value = "yPoint"
$ws.Range("J5").Value = "zPoint"
$ws.Range("B4").Value = "Delta"

# --- 5. Misc view state, matching the saved workbook ------------------------
$ws.Range("O13").Select()
